# Adds a new "2022-Q4" quarter sheet (with fund holding data) and records it
# as a new top row in the "总计" (totals) summary sheet, shifting the existing
# quarter rows/sheets down by one - per commit "feat: add 2022-Q4 data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while preserving text semantics for
# strings that look numeric (e.g. "2.88", "004317"), mirroring how the
# source workbook stores these figures as inline strings rather than
# numbers. Plain numbers are written as real numbers.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

function Set-NumValue($range, $value) {
    $range.Value = $value
}

# ===========================================================================
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 on top
#    of the existing quarterly rows.
# ===========================================================================
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

Set-NumValue $summary.Range("A2") 0
Set-TextValue $summary.Range("B2") "2022-Q4"
Set-NumValue $summary.Range("C2") 6
Set-NumValue $summary.Range("D2") 0.24

# Restore the bold/bordered look used on the rest of column A by copying the
# format from the row below (which already carries the correct style).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Column A is a simple 0-based row index; renumber the rows that shifted
# down to keep it consistent (1..6 for the former rows 2..7).
for ($r = 3; $r -le 8; $r++) {
    $idx = $r - 2
    Set-NumValue $summary.Range("A$r") $idx
}

# ===========================================================================
# 2. Insert a brand-new worksheet named "2022-Q4" right before the existing
#    "2022-Q2" sheet and populate it with the quarter's fund holding data.
#    The sheet is created by copying the "2022-Q2" sheet so that it inherits
#    the same layout/formatting (sheetPr, page margins, header/column
#    styles) already used by every other quarter sheet in this workbook.
# ===========================================================================
$existingSecond = $wb.Worksheets.Item(2)
$existingSecond.Copy($existingSecond)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Clear the copied data rows (2-6); the header row (1) and its formatting
# are reused as-is.
$q4.Range("A2:H6").ClearContents()

# -- Header row -------------------------------------------------------------
Set-TextValue $q4.Range("B1") "基金代码"
Set-TextValue $q4.Range("C1") "基金名称"
Set-TextValue $q4.Range("D1") "基金规模"
Set-TextValue $q4.Range("E1") "股票总仓位"
Set-TextValue $q4.Range("F1") "仓位占比"
Set-TextValue $q4.Range("G1") "持有市值(亿元)"
Set-TextValue $q4.Range("H1") "仓位排名"

# Re-apply the bold/bordered header style (ClearContents keeps formatting,
# but Set-TextValue resets the style of each cell it touches).
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# -- Data rows ----------------------------------------------------------
$rows = @(
    @("004317", "前海开源沪港深裕鑫灵活配置混合C", "2.88", "90.85", "3.60", "0.1037", 1),
    @("004316", "前海开源沪港深裕鑫灵活配置混合A", "2.30", "90.85", "3.60", "0.0828", 1),
    @("006923", "前海开源沪港深非周期性行业股票A", "0.28", "90.65", "7.63", "0.0214", 3),
    @("006924", "前海开源沪港深非周期性行业股票C", "0.24", "90.65", "7.63", "0.0183", 3),
    @("005493", "鑫元价值精选灵活配置混合A",       "0.57", "86.30", "1.89", "0.0108", 10),
    @("005494", "鑫元价值精选灵活配置混合C",       "0.00", "86.30", "1.89", "__NUM0__", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    Set-NumValue $q4.Range("A$r") $i

    Set-TextValue $q4.Range("B$r") $row[0]
    Set-TextValue $q4.Range("C$r") $row[1]
    Set-TextValue $q4.Range("D$r") $row[2]
    Set-TextValue $q4.Range("E$r") $row[3]
    Set-TextValue $q4.Range("F$r") $row[4]

    if ($row[5] -eq "__NUM0__") {
        Set-NumValue $q4.Range("G$r") 0
    } else {
        Set-TextValue $q4.Range("G$r") $row[5]
    }

    Set-NumValue $q4.Range("H$r") $row[6]
}

# Give column A its bold/bordered style to match the rest of the workbook
# (needed in particular for row 7, which didn't exist in the copied sheet).
$summary.Range("A3").Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)

# Restore the originally active tab ("总计", the first sheet) since creating
#/copying sheets above shifted the active tab to the newly added sheet.
$wb.Worksheets.Item(1).Activate()

